# Generate Report for Handback
# The "75af15d5-9871-42fc-9627-421b87f9cd98" file finished handback (in sync
# with en-US) for both zh-cn and de-de locales. Update the Overview sheet and
# both locale report sheets accordingly:
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - Latest Handback DateTime gets a fresh timestamp for each locale
#   - The stale "handback file not latest" Error Detail is cleared

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-18 20:50:12"
$zhcn.Range("P3").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-18 20:50:23"
$dede.Range("P3").Value = ""

$zhcn.Columns.Item(16).AutoFit() | Out-Null
$dede.Columns.Item(16).AutoFit() | Out-Null

Write-Output "Generate Report for Handback: updated handback status/timestamps"
